$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "From Date" (C) and "End Date" (D) values - rows moved one calendar
# month forward (Jun -> Jul etc.) along with the recalculated blank-day ranges.
$ws.Cells.Item(2, 3).Value2 = 45839
$ws.Cells.Item(2, 4).Value2 = 45991

$ws.Cells.Item(3, 3).Value2 = 45809
$ws.Cells.Item(3, 4).Value2 = 45991

$ws.Cells.Item(4, 3).Value2 = 45809
$ws.Cells.Item(4, 4).Value2 = 45991

$ws.Cells.Item(5, 3).Value2 = 45839
$ws.Cells.Item(5, 4).Value2 = 45991

$ws.Cells.Item(6, 3).Value2 = 45839
$ws.Cells.Item(6, 4).Value2 = 45991

$ws.Cells.Item(7, 3).Value2 = 45839
$ws.Cells.Item(7, 4).Value2 = 45930

$ws.Cells.Item(8, 3).Value2 = 45839
$ws.Cells.Item(8, 4).Value2 = 45869

$ws.Cells.Item(9, 3).Value2 = 45839
$ws.Cells.Item(9, 4).Value2 = 45869

$ws.Cells.Item(10, 3).Value2 = 45839
$ws.Cells.Item(10, 4).Value2 = 45930

$ws.Cells.Item(11, 3).Value2 = 45839
$ws.Cells.Item(11, 4).Value2 = 45930

$ws.Cells.Item(12, 3).Value2 = 45839
$ws.Cells.Item(12, 4).Value2 = 45930

$ws.Cells.Item(13, 3).Value2 = 45839
$ws.Cells.Item(13, 4).Value2 = 45869

$ws.Cells.Item(14, 3).Value2 = 45839
$ws.Cells.Item(14, 4).Value2 = 45869

# Update the selection on the sheet to reflect column D being selected.
$ws.Range("D1:D1048576").Select()
